$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44574
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7500
$ws.Range("S2").Value = 3750

$ws.Range("D3").Value = 44574
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44944
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7500
$ws.Range("S4").Value = 3750

$ws.Range("D5").Value = 44216
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 3500
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 3750
$ws.Range("S5").Value = 1875

$ws.Range("D6").Value = 44216
$ws.Range("L6").Value = "Segunda"
$ws.Range("N6").Value = 3000
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 3000
$ws.Range("S6").Value = 1500

$ws.Range("D7").Value = 44609
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 6500
$ws.Range("P7").Value = 6750
$ws.Range("S7").Value = 3375

$ws.Range("D8").Value = 44609
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("S8").Value = 3000

$ws.Range("D9").Value = 44532
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 5000

$ws.Range("D10").Value = 44532
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("S10").Value = 4000

$ws.Range("D11").Value = 44988
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 6500
$ws.Range("S11").Value = 3250

$ws.Range("D12").Value = 44195
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("O12").Value = 3500
$ws.Range("P12").Value = 3250
$ws.Range("S12").Value = 1625

$ws.Range("D13").Value = 44195
$ws.Range("L13").Value = "Segunda"
$ws.Range("N13").Value = 2500
$ws.Range("O13").Value = 2500
$ws.Range("P13").Value = 2500
$ws.Range("S13").Value = 1250

$ws.Range("D14").Value = 44602

$ws.Range("D15").Value = 44602
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 5000
$ws.Range("O15").Value = 5000
$ws.Range("P15").Value = 5000
$ws.Range("S15").Value = 2500

$ws.Range("D17").Value = 44917
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 7500
$ws.Range("P17").Value = 7250
$ws.Range("S17").Value = 3625

$ws.Range("D18").Value = 45014
$ws.Range("N18").Value = 7000
$ws.Range("O18").Value = 7500
$ws.Range("P18").Value = 7250
$ws.Range("S18").Value = 3625

$ws.Range("D19").Value = 44910
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 7500
$ws.Range("P19").Value = 7750
$ws.Range("S19").Value = 3875

$ws.Range("D20").Value = 44559
$ws.Range("N20").Value = 6000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 6500
$ws.Range("S20").Value = 3250

$ws.Range("D21").Value = 44559
$ws.Range("N21").Value = 5000
$ws.Range("O21").Value = 5000
$ws.Range("P21").Value = 5000
$ws.Range("S21").Value = 2500
